# Update cryptos list - GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force text storage so numeric-looking strings (e.g. "352.46") are not
    # silently coerced into floating point numbers by Excel's smart entry.
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "52.437.17"
$ws.Range("E2").Value = "  +1.92%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.922.84"
$ws.Range("E3").Value = "  +5.08%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.18%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "352.46"
$ws.Range("E5").Value = "  -0.13%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "112.64"
$ws.Range("E6").Value = "  +4.00%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.89%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.627"
$ws.Range("E9").Value = "  +0.96%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "40.25"
$ws.Range("E10").Value = "  +1.51%  "

# Row 11 - now Dogecoin (was TRON)
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D11") "0.0862"
$ws.Range("E11").Value = "  +3.49%  "

# Row 12 - now TRON (was Dogecoin)
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D12") "0.136"
$ws.Range("E12").Value = "  +0.79%  "

# Row 13 - Chainlink
Set-TextValue $ws.Range("D13") "20.19"
$ws.Range("E13").Value = "  +1.67%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "7.85"
$ws.Range("E14").Value = "  +1.83%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.377.70"
$ws.Range("E15").Value = "  +4.79%  "

# Row 16 - now Polygon (was WrappedEther)
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D16") "0.996"
$ws.Range("E16").Value = "  +6.39%  "

# Row 17 - now WrappedEther (was Polygon)
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.903.09"
$ws.Range("E17").Value = "  +3.82%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "52.438.39"
$ws.Range("E18").Value = "  +1.93%  "

# Row 19 - now InternetComputer(DFINITY) (was Uniswap)
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D19") "14.69"
$ws.Range("E19").Value = "  +9.84%  "

# Row 20 - ImmutableX
Set-TextValue $ws.Range("D20") "3.36"
$ws.Range("E20").Value = "  +6.88%  "

# Row 21 - now Uniswap (was InternetComputer(DFINITY))
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D21") "7.73"
$ws.Range("E21").Value = "  +0.84%  "

# Row 22 - ShibaInu
$ws.Range("D22").Value = "0.0₃0984"
$ws.Range("E22").Value = "  +1.75%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "71.16"
$ws.Range("E23").Value = "  +1.42%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "271.80"
$ws.Range("E24").Value = "  +2.13%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +2.08%  "

# Row 26 - EthereumClassic
Set-TextValue $ws.Range("D26") "27.01"
$ws.Range("E26").Value = "  +4.20%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.01%  "

# Row 28 - Kaspa
Set-TextValue $ws.Range("D28") "0.165"
$ws.Range("E28").Value = "  +0.29%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +3.22%  "

# Row 30 - InjectiveProtocol
Set-TextValue $ws.Range("D30") "38.08"
$ws.Range("E30").Value = "  +3.96%  "

# Row 31 - Toncoin
Set-TextValue $ws.Range("D31") "2.26"
$ws.Range("E31").Value = "  +1.54%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("D32") "6.46"
$ws.Range("E32").Value = "  +4.87%  "

# Row 33 - RenderToken
Set-TextValue $ws.Range("D33") "6.17"
$ws.Range("E33").Value = "  +9.89%  "

# Row 34 - OKB
Set-TextValue $ws.Range("D34") "53.13"
$ws.Range("E34").Value = "  +2.36%  "

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.0937"
$ws.Range("E35").Value = "  +9.77%  "

# Row 36 - VeChain
$ws.Range("E36").Value = "  +3.45%  "

# Row 37 - FirstDigitalUSD
Set-TextValue $ws.Range("D37") "0.998"
$ws.Range("E37").Value = "  -0.32%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  +7.24%  "

# Row 39 - Celestia
Set-TextValue $ws.Range("D39") "18.90"
$ws.Range("E39").Value = "  +0.97%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  +5.43%  "

# Row 41 - Stacks
Set-TextValue $ws.Range("D41") "2.74"
$ws.Range("E41").Value = "  +10.92%  "

# Row 42 - EnergySwap
Set-TextValue $ws.Range("D42") "24.30"
$ws.Range("E42").Value = "  +12.60%  "

# Row 43 - Stellar
$ws.Range("E43").Value = "  +2.44%  "

# Row 44 - Monero
Set-TextValue $ws.Range("D44") "122.87"
$ws.Range("E44").Value = "  +3.26%  "

# Row 45 - WEMIXToken
$ws.Range("E45").Value = "  +0.80%  "

# Row 46 - now NEARProtocol (was Maker)
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D46") "3.58"
$ws.Range("E46").Value = "  +6.36%  "

# Row 47 - now Maker (was NEARProtocol); E47 keeps its existing value (+5.08%)
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.220.77"

# Row 48 - ApeXProtocol
$ws.Range("E48").Value = "  +6.41%  "

# Row 49 - TheGraph
Set-TextValue $ws.Range("D49") "0.269"
$ws.Range("E49").Value = "  +25.87%  "

# Row 50 - SEI
Set-TextValue $ws.Range("D50") "0.957"
$ws.Range("E50").Value = "  +6.12%  "

# Row 51 - BEAM
Set-TextValue $ws.Range("D51") "0.0332"
$ws.Range("E51").Value = "  +15.04%  "
